$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '59.967.00'
$ws.Range("E2").Value = '  +0.25%  '
$ws.Range("D3").Value = '2.405.57'
$ws.Range("E3").Value = '  -0.47%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '553.69'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.47%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '135.67'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.98%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.583'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -0.87%  '
$ws.Range("E9").Value = '  -0.39%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '5.61'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -2.01%  '
$ws.Range("E12").Value = '  -1.57%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '24.61'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -0.39%  '
$ws.Range("D14").Value = '2.835.43'
$ws.Range("E14").Value = '  -0.41%  '
$ws.Range("D15").Value = '59.835.50'
$ws.Range("E16").Value = '  +0.23%  '
$ws.Range("D17").Value = '2.359.46'
$ws.Range("E17").Value = '  -2.23%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.16'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -1.10%  '
$ws.Range("E19").Value = '  +3.16%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '325.48'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -1.27%  '
$ws.Range("E21").Value = '  +1.07%  '
$ws.Range("E22").Value = '  -0.11%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '64.62'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -1.69%  '
$ws.Range("E24").Value = '  +4.39%  '
$ws.Range("E25").Value = '  +0.12%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -0.05%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.40'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +4.39%  '
$ws.Range("E28").Value = '  +0.92%  '
$ws.Range("E29").Value = '  -1.50%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '169.88'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -0.30%  '
$ws.Range("E31").Value = '  -0.86%  '
$ws.Range("E32").Value = '  +7.61%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.399'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -2.49%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '18.39'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -1.13%  '
$ws.Range("E35").Value = '  +0.02%  '
$ws.Range("E36").Value = '  +2.38%  '
$ws.Range("E37").Value = '  +0.03%  '
$ws.Range("E38").Value = '  +0.16%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '322.37'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +2.48%  '
$ws.Range("E40").Value = '  -0.68%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '146.54'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +5.97%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.58'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -1.77%  '
$ws.Range("E43").Value = '  -0.59%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '19.82'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +2.18%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0513'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -0.76%  '
$ws.Range("E46").Value = '  -0.75%  '
$ws.Range("E47").Value = '  -1.35%  '
$ws.Range("E48").Value = '  -0.05%  '
$ws.Range("E49").Value = '  -1.31%  '
$ws.Range("E50").Value = '  -0.66%  '
$ws.Range("E51").Value = '  -1.08%  '
